$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new price records were added to the daily Kiwi series (Macroferia
# Regional de Talca). This pushes the existing rows 201-283 down to
# rows 203-285, and the two new rows are populated below.
$ws.Rows("201:202").Insert()

# New row 201
$ws.Cells.Item(201, 1).Value = 5
$ws.Cells.Item(201, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(201, 3).Value = "Maule"
$ws.Cells.Item(201, 4).Value = 44755
$ws.Cells.Item(201, 5).Value = 7
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100101
$ws.Cells.Item(201, 8).Value = "Berries"
$ws.Cells.Item(201, 9).Value = 100101007
$ws.Cells.Item(201, 10).Value = "Kiwi"
$ws.Cells.Item(201, 11).Value = "Hayward"
$ws.Cells.Item(201, 12).Value = "Especial"
$ws.Cells.Item(201, 13).Value = 300
$ws.Cells.Item(201, 14).Value = 8000
$ws.Cells.Item(201, 15).Value = 8000
$ws.Cells.Item(201, 16).Value = 8000
$ws.Cells.Item(201, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(201, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(201, 19).Value = 444
$ws.Cells.Item(201, 20).Value = 18

# New row 202
$ws.Cells.Item(202, 1).Value = 5
$ws.Cells.Item(202, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(202, 3).Value = "Maule"
$ws.Cells.Item(202, 4).Value = 44755
$ws.Cells.Item(202, 5).Value = 7
$ws.Cells.Item(202, 6).Value = "Fruta"
$ws.Cells.Item(202, 7).Value = 100101
$ws.Cells.Item(202, 8).Value = "Berries"
$ws.Cells.Item(202, 9).Value = 100101007
$ws.Cells.Item(202, 10).Value = "Kiwi"
$ws.Cells.Item(202, 11).Value = "Hayward"
$ws.Cells.Item(202, 12).Value = "Primera"
$ws.Cells.Item(202, 13).Value = 200
$ws.Cells.Item(202, 14).Value = 6000
$ws.Cells.Item(202, 15).Value = 6000
$ws.Cells.Item(202, 16).Value = 6000
$ws.Cells.Item(202, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(202, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(202, 19).Value = 333
$ws.Cells.Item(202, 20).Value = 18
